$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("project hours")

# New "users" header in column E, styled like the other headers (B1:D1)
$ws.Range("E1").Value = "users"
$ws.Range("B1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# New user-list values for each project row
$ws.Range("E2").Value = "['Arun Lakshmanan', 'Hamid Jafarnejadsani', 'Sebastian Rodriguez']"
$ws.Range("E3").Value = "['Gavin Ananda']"
